$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C2").Value = 1.36
$ws.Range("C3").Value = 0.65
$ws.Range("C4").Value = 0.94
$ws.Range("C5").Value = 0.79
$ws.Range("C6").Value = 0.97
$ws.Range("C7").Value = 1.05
$ws.Range("C8").Value = 1.02
$ws.Range("C9").Value = 0.92
$ws.Range("C10").Value = 1.13
$ws.Range("C11").Value = 0.83
$ws.Range("C12").Value = 0.87
$ws.Range("C13").Value = 0.82

[void]$ws.Range("C14").Select()
